$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "파이썬 모든 패키지,라이브러리 업데이트 하기,Updating All Packages In A Virtual Environment on Windows"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2952"

$ws.Range("D36").Value = "ACGAN : Auxiliary Classifier GANs"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/353"

$ws.Range("D42").Value = "CUDA 설치 및 설치 실패 대처"
$ws.Range("E42").Value = "https://kjk92.tistory.com/82"
